$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J10").Value = 3
$ws.Range("M10").Value = '82.95 MPH'
$ws.Range("M12").Value = '14.74°'
$ws.Range("J14").Value = 'Roblez'
$ws.Range("M14").Value = 'Line Drive'
$ws.Range("J15").Value = 'Right'
$ws.Range("M15").Value = 'Single'
$ws.Range("J16").Value = '88-90 MPH'
$ws.Range("J17").Value = 'CB,FB,CH'
$ws.Range("J19").Value = 1
$ws.Range("M19").Value = '70.73 MPH'
$ws.Range("J20").Value = 1
$ws.Range("M21").Value = '-43.16°'
$ws.Range("J23").Value = 'Roblez'
$ws.Range("M23").Value = 'Ground Ball'
$ws.Range("M24").Value = 'Single'
$ws.Range("J25").Value = '88-90 MPH'
$ws.Range("J26").Value = 'CB,FB,CH'
$ws.Range("J28").Value = 9
$ws.Range("M28").Value = 'nan MPH'
$ws.Range("M30").Value = 'nan°'
$ws.Range("J32").Value = 'Thompson'
$ws.Range("M32").Value = 'Undefined'
$ws.Range("J33").Value = 'Left'
$ws.Range("M33").Value = 'Undefined'
$ws.Range("J34").Value = '84-84 MPH'
$ws.Range("J35").Value = 'SL,FB,CH'
$ws.Range("J37").Value = 8
$ws.Range("M37").Value = 'nan MPH'
$ws.Range("J38").Value = 0
$ws.Range("M39").Value = 'nan°'
$ws.Range("J41").Value = 'Thompson'
$ws.Range("M41").Value = 'Undefined'
$ws.Range("J42").Value = 'Left'
$ws.Range("M42").Value = 'Undefined'
$ws.Range("J43").Value = '84-84 MPH'
$ws.Range("J44").Value = 'SL,FB,CH'
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = '86.78 MPH'
$ws.Range("M48").Value = '22.04°'
$ws.Range("J50").Value = 'Plum'
$ws.Range("M50").Value = 'Line Drive'
$ws.Range("J51").Value = 'Right'
$ws.Range("M51").Value = 'Out'
$ws.Range("J52").Value = '84-86 MPH'
$ws.Range("J53").Value = 'SL,FB,CH'
$ws.Range("J61").Value = 5
$ws.Range("M61").Value = '83.46 MPH'
$ws.Range("M63").Value = '69.2°'
$ws.Range("J65").Value = 'Herbst'
$ws.Range("M65").Value = 'Popup'
$ws.Range("M66").Value = 'Out'
$ws.Range("J67").Value = '83-85 MPH'
$ws.Range("J68").Value = 'SL,CB,FB,CH'
$ws.Range("J70").Value = 4
$ws.Range("M70").Value = '93.65 MPH'
$ws.Range("M72").Value = '33.31°'
$ws.Range("M74").Value = 'Fly Ball'
$ws.Range("J76").Value = '88-90 MPH'
$ws.Range("J77").Value = 'CB,FB,CH'
